# Update "want to go" counts (column F) across the 展览 / 本地生活 / 全部类型 sheets,
# plus the special case on 展览!17 (CICF x AGF event) whose lowest price became
# sellable (column G) and whose cover image (column I) was refreshed.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibitions) sheet ---
$wsExhibit.Range("F6").Value  = 957
$wsExhibit.Range("F7").Value  = 482
$wsExhibit.Range("F9").Value  = 2254
$wsExhibit.Range("F10").Value = 641
$wsExhibit.Range("F13").Value = 1156
$wsExhibit.Range("F15").Value = 2269
$wsExhibit.Range("F16").Value = 724

# Row 17: 广州·2024 CICF×AGF动漫游戏盛典 ...
$wsExhibit.Range("F17").Value = 15268
$wsExhibit.Range("G17").Value = 98
$wsExhibit.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202409/LXvfupUl1726284714478.jpeg"

$wsExhibit.Range("F19").Value = 1430
$wsExhibit.Range("F20").Value = 580
$wsExhibit.Range("F21").Value = 574
$wsExhibit.Range("F22").Value = 151
$wsExhibit.Range("F23").Value = 543
$wsExhibit.Range("F24").Value = 156
$wsExhibit.Range("F25").Value = 115
$wsExhibit.Range("F29").Value = 27
$wsExhibit.Range("F31").Value = 42

# --- 本地生活 (Local life) sheet ---
$wsLocal.Range("F2").Value = 5751
$wsLocal.Range("F3").Value = 502
$wsLocal.Range("F4").Value = 484

# --- 全部类型 (All types) sheet ---
$wsAll.Range("F3").Value  = 502
$wsAll.Range("F4").Value  = 484
$wsAll.Range("F7").Value  = 957
$wsAll.Range("F9").Value  = 482
$wsAll.Range("F11").Value = 2254
$wsAll.Range("F12").Value = 641
$wsAll.Range("F17").Value = 1156
$wsAll.Range("F22").Value = 2269
$wsAll.Range("F23").Value = 724
$wsAll.Range("F27").Value = 1430
$wsAll.Range("F28").Value = 580
$wsAll.Range("F29").Value = 574
$wsAll.Range("F30").Value = 151
$wsAll.Range("F31").Value = 543
$wsAll.Range("F32").Value = 156
$wsAll.Range("F33").Value = 115
$wsAll.Range("F41").Value = 27
$wsAll.Range("F49").Value = 42
